$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows appended to the training log (rows 60-64)
$ws.Range("A60").Value = [DateTime]::ParseExact("2023-05-06", "yyyy-MM-dd", $null)
$ws.Range("B60").Value = "Your First Solo"
$ws.Range("C60").Value = "Collision Avoidance"

$ws.Range("A61").Value = [DateTime]::ParseExact("2023-08-30", "yyyy-MM-dd", $null)
$ws.Range("B61").Value = "Your First Solo"
$ws.Range("C61").Value = "Student Pilot & Medical Cetificate"

$ws.Range("A62").Value = [DateTime]::ParseExact("2023-08-30", "yyyy-MM-dd", $null)
$ws.Range("B62").Value = "Your First Solo"
$ws.Range("C62").Value = "Air Facts: Fit for Flight"

$ws.Range("A63").Value = [DateTime]::ParseExact("2023-08-30", "yyyy-MM-dd", $null)
$ws.Range("B63").Value = "Your First Solo"
$ws.Range("C63").Value = "Solo"

$ws.Range("A64").Value = [DateTime]::ParseExact("2023-08-31", "yyyy-MM-dd", $null)
$ws.Range("B64").Value = "Your Dual Cross Countries"

# Match the source date format used throughout column A (same built-in
# date style already used by the existing rows, s="1")
$ws.Range("A60:A64").NumberFormat = "m/d/yy"

# Column B got a bit wider to fit the new module name (stored width ends up
# at exactly 24 once the host's char->pixel-grid rounding is applied)
$ws.Columns.Item(2).ColumnWidth = 23.17

# Refresh the frozen pane view to the new scroll position and selection
$ws.Range("C65").Select()
$ws.Application.ActiveWindow.ScrollRow = 38
